$wb = $excel.ActiveWorkbook

# --- Characters sheet: add new "sight" stat column (K) ---
$ws1 = $wb.Worksheets.Item("Characters")

$ws1.Range("K1").Value = "sight"
$ws1.Range("K2").Value = 3
$ws1.Range("K3").Value = 3
$ws1.Range("K4").Value = 3

# --- View state: Characters becomes the active/selected sheet & cell ---
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws1.Range("K8").Select() | Out-Null
